$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: PAREJAS LIBRE  (fix team-name label typo) - edited first so that the
# workbook ends up with TRIO MASCULINO as the active/selected tab, as before.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("PAREJAS LIBRE")
$ws2.Range("C5").Value = "FEEDBACK LATINO"

$ws2.Activate()
$ws2.Range("C6").Select()

# ---------------------------------------------------------------------------
# Sheet: TRIO MASCULINO  (scores refreshed + sorted by total descending)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("TRIO MASCULINO")

# Row 3 - PASION LATINA (was FEEDBACK LATINO)
$ws1.Range("B3").Value = "PASION LATINA"
$ws1.Range("F3").Value = 10
$ws1.Range("G3").Value = 8
$ws1.Range("H3").Value = 8
$ws1.Range("I3").Value = 9
$ws1.Range("J3").Value = 9
$ws1.Range("K3").Value = 7
$ws1.Range("L3").Value = 8

# Row 4 - S.C LA VICTORIA (was SON LATINO ZULIA)
$ws1.Range("B4").Value = "S.C LA VICTORIA"
$ws1.Range("F4").Value = 10
$ws1.Range("G4").Value = 8
$ws1.Range("H4").Value = 8
$ws1.Range("I4").Value = 9
$ws1.Range("J4").Value = 9
$ws1.Range("K4").Value = 6
$ws1.Range("L4").Value = 8

# Row 5 - FEEDBACK LATINO (was S.C LA VICTORIA)
$ws1.Range("B5").Value = "FEEDBACK LATINO"
$ws1.Range("F5").Value = 7
$ws1.Range("G5").Value = 7
$ws1.Range("H5").Value = 7
$ws1.Range("I5").Value = 8
$ws1.Range("J5").Value = 8
$ws1.Range("K5").Value = 6
$ws1.Range("L5").Value = 7

# Row 6 - SON LATINO ZULIA (was PASION LATINA)
$ws1.Range("B6").Value = "SON LATINO ZULIA"
$ws1.Range("F6").Value = 7
$ws1.Range("G6").Value = 6
$ws1.Range("H6").Value = 7
$ws1.Range("I6").Value = 7
$ws1.Range("J6").Value = 7
$ws1.Range("K6").Value = 6
$ws1.Range("L6").Value = 6

# Record the sort (data is already placed in the sorted order above) so the
# worksheet keeps a <sortState> the same way Excel does after Data > Sort.
$sort1 = $ws1.Sort
$sort1.SortFields.Clear()
$sort1.SortFields.Add($ws1.Range("D3:D6"), 0, 2) | Out-Null
$sort1.SetRange($ws1.Range("B3:L6"))
$sort1.Header = 2
$sort1.Apply()

$ws1.Activate()
$ws1.Range("B2:L6").Select()
